$wb = $excel.ActiveWorkbook

# --- Sheet "User Stories Planejadas": update Sprint values (col E) ---
$wsPlan = $wb.Worksheets.Item("User Stories Planejadas")
$wsPlan.Range("E5").Value = 3
$wsPlan.Range("E6").Value = 3
$wsPlan.Range("E7").Value = 4
$wsPlan.Range("E8").Value = 4
$wsPlan.Range("E9").Value = 5

# Update the stored cursor/selection for this sheet (without activating it -
# it is not the active tab before or after the edit).
$wsPlan.Range("B41").Select()

# --- Sheet "User Stories Realizadas": clear the two stray trailing rows ---
$wsReal = $wb.Worksheets.Item("User Stories Realizadas")
$wsReal.Range("A10:A11").ClearContents()

# Make this sheet the active tab and set its selection/scroll position.
$wsReal.Activate()
$wsReal.Range("B8").Select()
